$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
try {
  $ws.StandardWidth = 9
  Write-Output "ok StandardWidth"
} catch {
  Write-Output "ERR StandardWidth: $_"
}
